$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-27 Thursday" "2025-03-28 Friday"

Replace-Text "602×5=3010" "948×9=8532"
Replace-Text "122×7=854" "246×7=1722"
Replace-Text "164×9=1476" "848×3=2544"
Replace-Text "527×2=1054" "755×6=4530"
Replace-Text "477×2=954" "490×4=1960"

Replace-Text "205×9=1845" "970×7=6790"
Replace-Text "939×7=6573" "898×6=5388"
Replace-Text "518×7=3626" "686×6=4116"
Replace-Text "629×6=3774" "748×6=4488"
Replace-Text "784×4=3136" "817×3=2451"

Replace-Text "953×5=4765" "620×7=4340"
Replace-Text "833×2=1666" "249×3=747"
Replace-Text "650×2=1300" "480×3=1440"
Replace-Text "381×7=2667" "577×2=1154"
Replace-Text "798×3=2394" "769×3=2307"

Replace-Text "547×2=1094" "408×3=1224"
Replace-Text "319×6=1914" "864×3=2592"
Replace-Text "411×8=3288" "892×5=4460"
Replace-Text "664×3=1992" "135×7=945"
Replace-Text "469×6=2814" "203×3=609"

Replace-Text "363×6=2178" "990×7=6930"
Replace-Text "248×4=992" "121×7=847"
Replace-Text "862×5=4310" "444×8=3552"
Replace-Text "798×6=4788" "610×4=2440"
Replace-Text "715×9=6435" "668×5=3340"
